$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metrics sheet: update the underlying metric values (B2:B13) and move the
# lingering selection from D31 to D12.
# ---------------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsMetrics.Range("B2").Value = 67070.19
$wsMetrics.Range("B3").Value = 55337.98
$wsMetrics.Range("B4").Value = 21081.43
$wsMetrics.Range("B5").Value = 2630
$wsMetrics.Range("B6").Value = 4434201.66
$wsMetrics.Range("B7").Value = 3745156.65
$wsMetrics.Range("B8").Value = 1291683.57
$wsMetrics.Range("B9").Value = 171631
$wsMetrics.Range("B10").Value = 32899525.460999828
$wsMetrics.Range("B11").Value = 31020378.170000002
$wsMetrics.Range("B12").Value = 11573392.460000001
$wsMetrics.Range("B13").Value = 1269258

# ---------------------------------------------------------------------------
# today sheet: the "month" figures (B11:B14) and their helper columns
# (E11:E14) were wiped out, so clear them - the dependent F column formulas
# stay in place and simply recompute to 0. The rest of the sheet (B15:B22,
# E15:E22, F15:F22) recomputes automatically once Metrics changes ripple
# through the =Metrics!Bx references.
# ---------------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Range("B11:B14").ClearContents()
$wsToday.Range("E11:E14").ClearContents()

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping: the saved selections moved, and the
# active tab moved from "today" to "bksr".
# ---------------------------------------------------------------------------
$wsMetrics.Range("D12").Select()
$wsToday.Range("F7").Select()

$wsBksr = $wb.Worksheets.Item("bksr")
$wsBksr.Activate()
$wsBksr.Range("G32").Select()
